$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.634.36'
$ws.Range("E2").Value = '  +1.96%  '
$ws.Range("D3").Value = '1.807.24'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.75'
$ws.Range("E5").Value = '  -3.00%  '
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4351'
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3761'
$ws.Range("E8").Value = '  +7.12%  '
$ws.Range("E9").Value = '  -2.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07683'
$ws.Range("E10").Value = '  +3.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.139'
$ws.Range("E11").Value = '  -0.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.60'
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.275'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.558'
$ws.Range("E15").Value = '  +3.88%  '
$ws.Range("D16").Value = '1.810.09'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001094'
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06749'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.06'
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.63'
$ws.Range("E21").Value = '  +2.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.287'
$ws.Range("D23").Value = '28.638.60'
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.76'
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.448'
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.57'
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.07'
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.353'
$ws.Range("E28").Value = '  -5.37%  '
$ws.Range("D29").Value = '2.018.75'
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.271'
$ws.Range("E30").Value = '  -1.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.76'
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.971'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.816'
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09196'
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2229'
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.16'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06338'
$ws.Range("E37").Value = '  +1.60%  '
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6602'
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02317'
$ws.Range("E39").Value = '  -2.01%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.208'
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.200'
$ws.Range("E41").Value = '  -1.55%  '
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.439'
$ws.Range("E42").Value = '  -2.90%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.070'
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.04'
$ws.Range("E45").Value = '  -0.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6075'
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.800'
$ws.Range("E47").Value = '  -1.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.06'
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.031'
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07095'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.19'
$ws.Range("E51").Value = '  -0.14%  '